# Generate Report for handback
#
# The d27a22fc-be9e-49a6-b335-595d2912778f.md file has now been handed back
# (both zh-cn and de-de). Update the Overview sheet status for that row, and
# fill in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the per-locale report sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: row 6 is the d27a22fc-... file. It has now been handed
# back for both locales, so its status changes from "Ready for handoff" to
# "Handed back: in sync with en-US".
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B6").Value = "Handed back: in sync with en-US"
$overview.Range("C6").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# zh-cn sheet: row 6 (d27a22fc-be9e-49a6-b335-595d2912778f.md)
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("B6").Value = "Handed back: in sync with en-US"

$zhcn.Range("E6").Value = "d27a22fc-be9e-49a6-b335-595d2912778f.md"
$zhcn.Hyperlinks.Add(
    $zhcn.Range("E6"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a05e351bb4e5a49bb11e461e25e8e3fbb06e9e4f/e2e/d27a22fc-be9e-49a6-b335-595d2912778f.md",
    "",
    "",
    "d27a22fc-be9e-49a6-b335-595d2912778f.md"
) | Out-Null

$zhcn.Range("F6").Value = "d27a22fc-be9e-49a6-b335-595d2912778f.850905496532c297fdf8772bb3c5a73b9f2bba2d.zh-cn.xlf"
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F6"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f4f64f7e06c0f8a8c1a9d8b8ddcd16b0b6a85f3a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/d27a22fc-be9e-49a6-b335-595d2912778f.850905496532c297fdf8772bb3c5a73b9f2bba2d.zh-cn.xlf",
    "",
    "",
    "d27a22fc-be9e-49a6-b335-595d2912778f.850905496532c297fdf8772bb3c5a73b9f2bba2d.zh-cn.xlf"
) | Out-Null

$zhcn.Range("G6").Value = "2016-01-28 09:24:32"

# ---------------------------------------------------------------------------
# de-de sheet: row 6 (d27a22fc-be9e-49a6-b335-595d2912778f.md)
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("B6").Value = "Handed back: in sync with en-US"

$dede.Range("E6").Value = "d27a22fc-be9e-49a6-b335-595d2912778f.md"
$dede.Hyperlinks.Add(
    $dede.Range("E6"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/2b5e3f7a6f0d8c2c9f4b9a7ad5e6e31c6f9f2d8a/e2e/d27a22fc-be9e-49a6-b335-595d2912778f.md",
    "",
    "",
    "d27a22fc-be9e-49a6-b335-595d2912778f.md"
) | Out-Null

$dede.Range("F6").Value = "d27a22fc-be9e-49a6-b335-595d2912778f.850905496532c297fdf8772bb3c5a73b9f2bba2d.de-de.xlf"
$dede.Hyperlinks.Add(
    $dede.Range("F6"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9c3a7e2f5b6d4a8e1f0c9b7a6d5e4f3c2b1a0f9e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/d27a22fc-be9e-49a6-b335-595d2912778f.850905496532c297fdf8772bb3c5a73b9f2bba2d.de-de.xlf",
    "",
    "",
    "d27a22fc-be9e-49a6-b335-595d2912778f.850905496532c297fdf8772bb3c5a73b9f2bba2d.de-de.xlf"
) | Out-Null

$dede.Range("G6").Value = "2016-01-28 09:24:52"
